# Refresh the cryptocurrency Price (column D) and 1h Volume-change (column E)
# figures pulled from coinranking.com, as performed by the scheduled
# GitHub Actions workflow that updates this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is @(CellReference, NewTextValue).
$updates = @(
    @("D2", "26.885.12"),
    @("E2", "  -0.81%  "),
    @("D3", "1.861.77"),
    @("E3", "  -0.34%  "),
    @("E4", "  -0.15%  "),
    @("D5", "304.82"),
    @("E5", "  -0.89%  "),
    @("E6", "  -0.09%  "),
    @("D7", "0.5057"),
    @("E7", "  -0.12%  "),
    @("D8", "0.3627"),
    @("E8", "  -3.05%  "),
    @("D9", "0.07165"),
    @("E9", "  +0.17%  "),
    @("D10", "0.8955"),
    @("E10", "  +0.99%  "),
    @("E11", "  +0.38%  "),
    @("D12", "1.855.23"),
    @("E12", "  -0.69%  "),
    @("D13", "0.07470"),
    @("E13", "  -1.15%  "),
    @("D14", "92.62"),
    @("E14", "  +3.88%  "),
    @("E15", "  -1.64%  "),
    @("E16", "  -0.15%  "),
    @("D17", "0.000008465"),
    @("E17", "  -0.07%  "),
    @("D18", "14.14"),
    @("E18", "  +0.25%  "),
    @("D19", "0.9998"),
    @("E19", "  -0.14%  "),
    @("D20", "26.923.08"),
    @("E20", "  -0.87%  "),
    @("E21", "  -0.89%  "),
    @("D22", "2.091.78"),
    @("E22", "  -0.03%  "),
    @("E23", "  -1.95%  "),
    @("D24", "6.412"),
    @("E24", "  -1.06%  "),
    @("D25", "147.94"),
    @("E25", "  -1.76%  "),
    @("D26", "1.794"),
    @("E26", "  -2.39%  "),
    @("D27", "17.88"),
    @("E27", "  -0.36%  "),
    @("E28", "  -1.73%  "),
    @("D29", "113.07"),
    @("E29", "  +0.34%  "),
    @("D30", "4.676"),
    @("E30", "  -1.54%  "),
    @("D31", "4.677"),
    @("E31", "  -0.28%  "),
    @("D32", "0.09263"),
    @("E32", "  +2.38%  "),
    @("E33", "  -0.64%  "),
    @("D34", "2.981"),
    @("E34", "  -3.68%  "),
    @("D35", "0.7421"),
    @("E35", "  +0.49%  "),
    @("D36", "1.149"),
    @("E36", "  -0.87%  "),
    @("E37", "  +7.66%  "),
    @("D38", "0.01999"),
    @("E38", "  -1.74%  "),
    @("D39", "2.499"),
    @("E39", "  +0.12%  "),
    @("D40", "0.5552"),
    @("E40", "  +3.84%  "),
    @("E41", "  -0.69%  "),
    @("D42", "118.62"),
    @("E42", "  +2.56%  "),
    @("D43", "6.476"),
    @("E43", "  -1.79%  "),
    @("D44", "8.502"),
    @("E44", "  +2.24%  "),
    @("D45", "0.1467"),
    @("E45", "  -0.45%  "),
    @("D46", "0.4695"),
    @("E46", "  +1.36%  "),
    @("E47", "  -0.08%  "),
    @("D48", "10.07"),
    @("E48", "  +1.22%  "),
    @("D49", "1.563"),
    @("E49", "  -0.15%  "),
    @("D50", "36.97"),
    @("E50", "  +1.47%  "),
    @("D51", "62.98"),
    @("E51", "  -2.38%  ")
)

foreach ($update in $updates) {
    $cellRef = $update[0]
    $newValue = $update[1]
    $cell = $ws.Range($cellRef)
    # Force text storage (matching the existing inline-string cells) so
    # values like "304.82" or "26.885.12" are not reinterpreted by Excel
    # as numbers, losing precision/trailing zeros or using exponent form.
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.Style = "Normal"
}
